$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The document gets a big new "dziennik" entry appended at the end:
# a brand new final paragraph with the "OSTATECZNIE ..." text, and the
# hidden "_GoBack" bookmark (which currently sits collapsed right
# after "...programować." at the end of the very first paragraph)
# moves to sit collapsed right after the new final paragraph's text.
# ------------------------------------------------------------------

$newParaText = "OSTATECZNIE metodą prób i błędów udało mi się naprawić ten problem. Po pierwsze – nadałem złą nazwę bazie w funkcji „Baza”, następnie okazało się, że w funkcji która określa jaką kwerendę wykonuje na bazie danych były złe pola, oraz jeszcze w funkcji do pól w bazie danych pojawiła się gafa, ponownie nazwy pól, które nie występują w mojej bazie danych."

# 1. Append a new paragraph after the current last paragraph
#    ("Czeka mnie czytanie poradnikow"). We add a one-character
#    placeholder "X" after the real text so that, while we are
#    adding the bookmark below, the target insertion point is not
#    literally the paragraph's end (a collapsed range placed exactly
#    at a paragraph boundary gets mis-resolved by this host's
#    Bookmarks.Add). We strip the placeholder right after.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last
$newLast.Range.Text = $newParaText + "X"

# 2. Remove the old _GoBack bookmark from its current location
#    (end of the first paragraph).
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# 3. Add the _GoBack bookmark back, collapsed, just before the
#    trailing placeholder character.
$newLast = $d.Paragraphs.Last
$pos = $newLast.Range.End - 2
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 4. Strip the one-character placeholder now that the bookmark is
#    safely anchored right after the real text.
$newLast = $d.Paragraphs.Last
$placeholderRange = $d.Range($newLast.Range.End - 2, $newLast.Range.End - 1)
$placeholderRange.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Final last paragraph: [" $d.Paragraphs.Last.Range.Text "]"
